$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data representing the ifoCAST full series evaluation (window shifted by one
# additional quarter: rows 2-10 take on the values previously held by the row
# below them, and row 11 receives the newly evaluated data point).
$data = @(
    @(-0.03027116410798613, 0.3513215654800215, 0.1490724330844395, 0.3860989938920322, 0.3994404986893751, 14),
    @(-0.04905665594513928, 0.2345188573562595, 0.07249634168634074, 0.269251446953105, 0.275555079267236, 13),
    @(-0.03952956703121076, 0.2973229505318943, 0.1197329025896884, 0.3460244248455424, 0.359044665623571, 12),
    @(-0.02815451101896177, 0.1902470767633485, 0.06574748725752844, 0.2564127283453932, 0.267301873999216, 11),
    @(-0.04837455801976177, 0.3441411292542904, 0.1571651242787223, 0.3964405684068197, 0.4147623501812882, 10),
    @(-0.0723091577469028, 0.3260979148759209, 0.116774734798478, 0.3417231844614556, 0.3542448102195189, 9),
    @(-0.04314485034007026, 0.2928395475983343, 0.1273075999664306, 0.3568019057774644, 0.3786383451347186, 8),
    @(-0.041375596389022, 0.1988059984120452, 0.06577567631596841, 0.2564676905888311, 0.2733880497193041, 7),
    @(-0.09920882285969852, 0.4033304192040505, 0.2097558308966996, 0.4579910816781257, 0.4897919236126613, 6),
    @(-0.009286038702706412, 0.2427162589609241, 0.07139578606315643, 0.2671998990702587, 0.2985581089762265, 5)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
